{"js": "// Replace the three-digit x one-digit multiplication equations in the\n// table with the new values from the commit. Each old value occurs\n// exactly once in the document, so a plain (non-wildcard) search +\n// whole-match replace is unambiguous for every pair.\nconst replacements = [\n  [\"943\u00d73=2829\", \"434\u00d79=3906\"],\n  [\"564\u00d74=2256\", \"675\u00d74=2700\"],\n  [\"897\u00d75=4485\", \"367\u00d79=3303\"],\n  [\"227\u00d73=681\", \"845\u00d75=4225\"],\n  [\"394\u00d74=1576\", \"867\u00d75=4335\"],\n  [\"268\u00d75=1340\", \"581\u00d77=4067\"],\n  [\"823\u00d74=3292\", \"747\u00d72=1494\"],\n  [\"976\u00d77=6832\", \"403\u00d73=1209\"],\n  [\"869\u00d74=3476\", \"406\u00d73=1218\"],\n  [\"551\u00d79=4959\", \"769\u00d79=6921\"],\n  [\"408\u00d77=2856\", \"253\u00d74=1012\"],\n  [\"218\u00d78=1744\", \"325\u00d72=650\"],\n  [\"293\u00d73=879\", \"782\u00d75=3910\"],\n  [\"456\u00d74=1824\", \"832\u00d78=6656\"],\n  [\"693\u00d79=6237\", \"857\u00d72=1714\"],\n  [\"911\u00d74=3644\", \"485\u00d73=1455\"],\n  [\"375\u00d75=1875\", \"502\u00d76=3012\"],\n  [\"573\u00d76=3438\", \"269\u00d77=1883\"],\n  [\"973\u00d76=5838\", \"460\u00d79=4140\"],\n  [\"604\u00d77=4228\", \"253\u00d74=1012\"],\n  [\"899\u00d74=3596\", \"594\u00d72=1188\"],\n  [\"831\u00d75=4155\", \"811\u00d79=7299\"],\n  [\"135\u00d78=1080\", \"921\u00d73=2763\"],\n  [\"252\u00d76=1512\", \"636\u00d75=3180\"],\n  [\"178\u00d76=1068\", \"554\u00d73=1662\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication equations in the\n# table with the new values from the commit. Each old value occurs\n# exactly once in the document, so Find/Replace (wdReplaceOne) on the\n# whole document Range is unambiguous for every pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"943\u00d73=2829\", \"434\u00d79=3906\"),\n    @(\"564\u00d74=2256\", \"675\u00d74=2700\"),\n    @(\"897\u00d75=4485\", \"367\u00d79=3303\"),\n    @(\"227\u00d73=681\", \"845\u00d75=4225\"),\n    @(\"394\u00d74=1576\", \"867\u00d75=4335\"),\n    @(\"268\u00d75=1340\", \"581\u00d77=4067\"),\n    @(\"823\u00d74=3292\", \"747\u00d72=1494\"),\n    @(\"976\u00d77=6832\", \"403\u00d73=1209\"),\n    @(\"869\u00d74=3476\", \"406\u00d73=1218\"),\n    @(\"551\u00d79=4959\", \"769\u00d79=6921\"),\n    @(\"408\u00d77=2856\", \"253\u00d74=1012\"),\n    @(\"218\u00d78=1744\", \"325\u00d72=650\"),\n    @(\"293\u00d73=879\", \"782\u00d75=3910\"),\n    @(\"456\u00d74=1824\", \"832\u00d78=6656\"),\n    @(\"693\u00d79=6237\", \"857\u00d72=1714\"),\n    @(\"911\u00d74=3644\", \"485\u00d73=1455\"),\n    @(\"375\u00d75=1875\", \"502\u00d76=3012\"),\n    @(\"573\u00d76=3438\", \"269\u00d77=1883\"),\n    @(\"973\u00d76=5838\", \"460\u00d79=4140\"),\n    @(\"604\u00d77=4228\", \"253\u00d74=1012\"),\n    @(\"899\u00d74=3596\", \"594\u00d72=1188\"),\n    @(\"831\u00d75=4155\", \"811\u00d79=7299\"),\n    @(\"135\u00d78=1080\", \"921\u00d73=2763\"),\n    @(\"252\u00d76=1512\", \"636\u00d75=3180\"),\n    @(\"178\u00d76=1068\", \"554\u00d73=1662\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
